$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Info")

# Wrap the @base / @prefix URL values (column D) in angle brackets.
$ws.Range("D1").Value = "<http://example.org/ex1>"
$ws.Range("D2").Value = "<http://foo.bar/data#>"
$ws.Range("D3").Value = "<http://foo.bar/model#>"
$ws.Range("D4").Value = "<http:/foo.bar/data#>"
$ws.Range("D5").Value = "<http://foo.bar/model#>"

# D1:D3 drop their right-aligned style (revert to the default/normal style),
# matching D4:D5 which were already unstyled.
$ws.Range("D1:D3").Style = "Normal"

# Widen column D to fit the longer, bracketed values.
$ws.Columns("D").ColumnWidth = 26.88

# Select the whole of column D, as recorded in the saved view state.
$ws.Columns("D").Select()
